$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Extract" label (Text Box 88) moves from (301.5, 337.125) pt to (318, 312) pt
$extract = $s.Shapes.Item("Text Box 88")
$extract.Left = 318
$extract.Top = 312

# "Re-index" label (Text Box 89) moves from (439.5, 334.25) pt to (444, 312) pt
$reindex = $s.Shapes.Item("Text Box 89")
$reindex.Left = 444
$reindex.Top = 312
